$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.540008666666668
$ws.Range("H2").Value = 19.620026
$ws.Range("I2").Value = 0.2365207520404831
$ws.Range("J2").Value = 0.2365207520404831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.425703666666666
$ws.Range("N2").Value = 4.277111
$ws.Range("O2").Value = 0.04715501820393346
$ws.Range("P2").Value = 0.04715501820393346
$ws.Range("Q2").Value = 9.324114336098445
$ws.Range("R2").Value = 83.917029024886
$ws.Range("S2").Value = 0.01115314036807701
$ws.Range("T2").Value = 0.01115314036807701

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.540008666666668
$ws.Range("H3").Value = 19.620026
$ws.Range("I3").Value = 0.2365207520404831
$ws.Range("J3").Value = 0.2365207520404831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.61433933333333
$ws.Range("N3").Value = 61.843018
$ws.Range("O3").Value = 0.6818173855147049
$ws.Range("P3").Value = 0.6818173855147048
$ws.Range("Q3").Value = 134.8179578976076
$ws.Range("R3").Value = 1213.361621078468
$ws.Range("S3").Value = 0.161263960776214
$ws.Range("T3").Value = 0.1612639607762139

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.540008666666668
$ws.Range("H4").Value = 19.620026
$ws.Range("I4").Value = 0.2365207520404831
$ws.Range("J4").Value = 0.2365207520404831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.194356666666666
$ws.Range("N4").Value = 24.58307
$ws.Range("O4").Value = 0.2710275962813615
$ws.Range("P4").Value = 0.2710275962813615
$ws.Range("Q4").Value = 53.59116361775779
$ws.Range("R4").Value = 482.3204725598201
$ws.Range("S4").Value = 0.06410365089619206
$ws.Range("T4").Value = 0.06410365089619206

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.67485766666667
$ws.Range("H5").Value = 44.024573
$ws.Range("I5").Value = 0.5307192311682535
$ws.Range("J5").Value = 0.5307192311682536
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.425703666666666
$ws.Range("N5").Value = 4.277111
$ws.Range("O5").Value = 0.04715501820393346
$ws.Range("P5").Value = 0.04715501820393346
$ws.Range("Q5").Value = 20.92199838317811
$ws.Range("R5").Value = 188.297985448603
$ws.Range("S5").Value = 0.02502607500691656
$ws.Range("T5").Value = 0.02502607500691657

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.67485766666667
$ws.Range("H6").Value = 44.024573
$ws.Range("I6").Value = 0.5307192311682535
$ws.Range("J6").Value = 0.5307192311682536
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.61433933333333
$ws.Range("N6").Value = 61.843018
$ws.Range("O6").Value = 0.6818173855147049
$ws.Range("P6").Value = 0.6818173855147048
$ws.Range("Q6").Value = 302.5124956090349
$ws.Range("R6").Value = 2722.612460481314
$ws.Range("S6").Value = 0.3618535986375129
$ws.Range("T6").Value = 0.3618535986375129

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.67485766666667
$ws.Range("H7").Value = 44.024573
$ws.Range("I7").Value = 0.5307192311682535
$ws.Range("J7").Value = 0.5307192311682536
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.194356666666666
$ws.Range("N7").Value = 24.58307
$ws.Range("O7").Value = 0.2710275962813615
$ws.Range("P7").Value = 0.2710275962813615
$ws.Range("Q7").Value = 120.2510177532344
$ws.Range("R7").Value = 1082.25915977911
$ws.Range("S7").Value = 0.143839557523824
$ws.Range("T7").Value = 0.143839557523824

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.436020999999999
$ws.Range("H8").Value = 19.308063
$ws.Range("I8").Value = 0.2327600167912634
$ws.Range("J8").Value = 0.2327600167912634
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.425703666666666
$ws.Range("N8").Value = 4.277111
$ws.Range("O8").Value = 0.04715501820393346
$ws.Range("P8").Value = 0.04715501820393346
$ws.Range("Q8").Value = 9.175858738443665
$ws.Range("R8").Value = 82.58272864599299
$ws.Range("S8").Value = 0.01097580282893988
$ws.Range("T8").Value = 0.01097580282893988

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.436020999999999
$ws.Range("H9").Value = 19.308063
$ws.Range("I9").Value = 0.2327600167912634
$ws.Range("J9").Value = 0.2327600167912634
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.61433933333333
$ws.Range("N9").Value = 61.843018
$ws.Range("O9").Value = 0.6818173855147049
$ws.Range("P9").Value = 0.6818173855147048
$ws.Range("Q9").Value = 132.6743208504593
$ws.Range("R9").Value = 1194.068887654134
$ws.Range("S9").Value = 0.158699826100978
$ws.Range("T9").Value = 0.158699826100978

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.436020999999999
$ws.Range("H10").Value = 19.308063
$ws.Range("I10").Value = 0.2327600167912634
$ws.Range("J10").Value = 0.2327600167912634
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.194356666666666
$ws.Range("N10").Value = 24.58307
$ws.Range("O10").Value = 0.2710275962813615
$ws.Range("P10").Value = 0.2710275962813615
$ws.Range("Q10").Value = 52.73905158815666
$ws.Range("R10").Value = 474.6514642934099
$ws.Range("S10").Value = 0.06308438786134546
$ws.Range("T10").Value = 0.06308438786134546
